$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("my_excel")

$ws.Range("D2").Value = 10.35
$ws.Range("E2").Value = 10.1
$ws.Range("D3").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
